$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card22")

# Row 9 updates: values move from "nan" placeholder to actual service entry
# D9 must stay a text value (not a number), so format as text, set it, then
# restore the default style so no extra formatting lingers on the cell.
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1107"
$ws.Range("D9").Style = "Normal"
$ws.Range("F9").Value = "Done "
$ws.Range("H9").Value = "Done "
$ws.Range("L9").Value = "15/2/2026"
$ws.Range("M9").Value = "سيرفيس"
$ws.Range("N9").Value = "تم تغيير الجرائد الاماميه وتغيير جريده 1خلفيه"
$ws.Range("O9").Value = "م محمد  ومحمود"

# Row 33 updates: previously empty cells now hold literal "nan" text
$ws.Range("B33").Value = "nan"
$ws.Range("C33").Value = "nan"
$ws.Range("D33").Value = "nan"
$ws.Range("E33").Value = "nan"
$ws.Range("F33").Value = "nan"
$ws.Range("G33").Value = "nan"
$ws.Range("H33").Value = "nan"
$ws.Range("I33").Value = "nan"
$ws.Range("J33").Value = "nan"
$ws.Range("K33").Value = "nan"
$ws.Range("P33").Value = "nan"
